$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "https://c123.com.br/bel-ar/FotoRetArq.asp?a=101051000%2Ejpg"
$ws.Range("B2").Value = 101051000

# Row 3
$ws.Range("A3").Value = "https://c123.com.br/bel-ar/FotoRetArq.asp?a=101021640%2Ejpg"
$ws.Range("B3").Value = 101021640

# Row 4
$ws.Range("A4").Value = "https://c123.com.br/bel-ar/FotoRetArq.asp?a=101041217%2Ejpg"
$ws.Range("B4").Value = 101041217

# Row 5
$ws.Range("A5").Value = "https://c123.com.br/bel-ar/FotoRetArq.asp?a=101042728%2Ejpg"
$ws.Range("B5").Value = 101042728

# Row 6
$ws.Range("A6").Value = "https://c123.com.br/bel-ar/FotoRetArq.asp?a=101021569%2Ejpg"
$ws.Range("B6").Value = 101021569

# Row 7
$ws.Range("A7").Value = "https://c123.com.br/bel-ar/FotoRetArq.asp?a=101021836%2Ejpg"
$ws.Range("B7").Value = 101021836

# Row 8
$ws.Range("A8").Value = "https://c123.com.br/bel-ar/FotoRetArq.asp?a=101043090%2Ejpg"
$ws.Range("B8").Value = 101043090

# Row 9
$ws.Range("A9").Value = "https://c123.com.br/bel-ar/FotoRetArq.asp?a=101021526%2Ejpg"
$ws.Range("B9").Value = 101021526

# Row 10
$ws.Range("A10").Value = "https://c123.com.br/bel-ar/FotoRetArq.asp?a=101023715%2Ejpg"
$ws.Range("B10").Value = 101023715

# Rows 11-12: clear previous leftover data (no more entries)
$ws.Range("A11:B12").ClearContents()

# Update selection to C6 as in the saved file
$ws.Range("C6").Select()
